$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "332.20"
Set-TextCell $ws.Range("E2") "0.39%"
Set-TextCell $ws.Range("G2") "13"

# Row 3
Set-TextCell $ws.Range("D3") "39.63"
Set-TextCell $ws.Range("E3") "-1.54%"
Set-TextCell $ws.Range("G3") "13"

# Row 4
Set-TextCell $ws.Range("D4") "5.767"
Set-TextCell $ws.Range("E4") "3.01%"
Set-TextCell $ws.Range("G4") "13"

# Row 5
Set-TextCell $ws.Range("D5") "0.08073"
Set-TextCell $ws.Range("E5") "-0.51%"
Set-TextCell $ws.Range("G5") "13"

# Row 6
Set-TextCell $ws.Range("D6") "2.013"
Set-TextCell $ws.Range("E6") "2.10%"
Set-TextCell $ws.Range("G6") "13"

# Row 7
Set-TextCell $ws.Range("B7") "GateToken"
Set-TextCell $ws.Range("C7") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell $ws.Range("D7") "4.519"
Set-TextCell $ws.Range("E7") "-0.78%"
Set-TextCell $ws.Range("G7") "13"

# Row 8
Set-TextCell $ws.Range("B8") "KuCoinToken"
Set-TextCell $ws.Range("C8") "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextCell $ws.Range("D8") "8.655"
Set-TextCell $ws.Range("E8") "-0.31%"
Set-TextCell $ws.Range("G8") "13"

# Row 9
Set-TextCell $ws.Range("D9") "2.953"
Set-TextCell $ws.Range("E9") "-0.81%"
Set-TextCell $ws.Range("G9") "13"

# Row 10
Set-TextCell $ws.Range("D10") "0.9228"
Set-TextCell $ws.Range("E10") "-2.73%"
Set-TextCell $ws.Range("G10") "13"

# Row 11
Set-TextCell $ws.Range("D11") "0.1288"
Set-TextCell $ws.Range("E11") "2.02%"
Set-TextCell $ws.Range("G11") "13"

# Row 12
Set-TextCell $ws.Range("D12") "0.1956"
Set-TextCell $ws.Range("E12") "-1.24%"
Set-TextCell $ws.Range("G12") "13"

# Row 13
Set-TextCell $ws.Range("D13") "8.764"
Set-TextCell $ws.Range("E13") "19.57%"
Set-TextCell $ws.Range("G13") "13"

# Row 14
Set-TextCell $ws.Range("D14") "0.09243"
Set-TextCell $ws.Range("E14") "0.50%"
Set-TextCell $ws.Range("G14") "13"

# Row 15
Set-TextCell $ws.Range("D15") "0.03547"
Set-TextCell $ws.Range("E15") "-1.43%"
Set-TextCell $ws.Range("G15") "13"

# Row 16
Set-TextCell $ws.Range("D16") "0.1053"
Set-TextCell $ws.Range("E16") "9.71%"
Set-TextCell $ws.Range("G16") "13"

# Row 17
Set-TextCell $ws.Range("D17") "0.001304"
Set-TextCell $ws.Range("E17") "-2.12%"
Set-TextCell $ws.Range("G17") "13"

# Row 18
Set-TextCell $ws.Range("D18") "0.006402"
Set-TextCell $ws.Range("E18") "2.09%"
Set-TextCell $ws.Range("G18") "13"

# Row 19
Set-TextCell $ws.Range("D19") "3.368"
Set-TextCell $ws.Range("E19") "0.14%"
Set-TextCell $ws.Range("G19") "13"

# Row 20
Set-TextCell $ws.Range("E20") "-1.00%"
Set-TextCell $ws.Range("G20") "13"

# Row 21
Set-TextCell $ws.Range("D21") "0.1361"
Set-TextCell $ws.Range("E21") "1.43%"
Set-TextCell $ws.Range("G21") "13"

# Row 22
Set-TextCell $ws.Range("D22") "0.2721"
Set-TextCell $ws.Range("E22") "11.02%"
Set-TextCell $ws.Range("G22") "13"

# Row 23
Set-TextCell $ws.Range("D23") "0.04420"
Set-TextCell $ws.Range("E23") "-0.23%"
Set-TextCell $ws.Range("G23") "13"

# Row 24
Set-TextCell $ws.Range("D24") "0.001261"
Set-TextCell $ws.Range("E24") "2.49%"
Set-TextCell $ws.Range("G24") "13"

# Row 25
Set-TextCell $ws.Range("D25") "0.004567"
Set-TextCell $ws.Range("E25") "5.37%"
Set-TextCell $ws.Range("G25") "13"

# Row 26
Set-TextCell $ws.Range("D26") "0.0001199"
Set-TextCell $ws.Range("E26") "-0.13%"
Set-TextCell $ws.Range("G26") "13"

# Row 27
Set-TextCell $ws.Range("G27") "13"

# Row 28
Set-TextCell $ws.Range("G28") "13"

# Row 29
Set-TextCell $ws.Range("G29") "13"

# Row 30
Set-TextCell $ws.Range("G30") "13"

# Row 31
Set-TextCell $ws.Range("G31") "13"

# Row 32
Set-TextCell $ws.Range("G32") "13"

# Row 33
Set-TextCell $ws.Range("G33") "13"

# Row 34
Set-TextCell $ws.Range("G34") "13"

# Row 35
Set-TextCell $ws.Range("G35") "13"

# Row 36
Set-TextCell $ws.Range("G36") "13"

# Row 37
Set-TextCell $ws.Range("G37") "13"

# Row 38
Set-TextCell $ws.Range("G38") "13"

# Row 39
Set-TextCell $ws.Range("D39") "0.02520"
Set-TextCell $ws.Range("E39") "0.18%"
Set-TextCell $ws.Range("G39") "13"

# Row 40
Set-TextCell $ws.Range("D40") "0.05474"
Set-TextCell $ws.Range("E40") "4.48%"
Set-TextCell $ws.Range("G40") "13"

# Row 41
Set-TextCell $ws.Range("D41") "0.007491"
Set-TextCell $ws.Range("E41") "-3.44%"
Set-TextCell $ws.Range("G41") "13"

# Row 42
Set-TextCell $ws.Range("D42") "0.009927"
Set-TextCell $ws.Range("E42") "10.65%"
Set-TextCell $ws.Range("G42") "13"

# Row 43
Set-TextCell $ws.Range("D43") "0.1411"
Set-TextCell $ws.Range("E43") "-1.51%"
Set-TextCell $ws.Range("G43") "13"

# Row 44
Set-TextCell $ws.Range("D44") "0.002106"
Set-TextCell $ws.Range("G44") "13"

# Row 45
Set-TextCell $ws.Range("D45") "0.01133"
Set-TextCell $ws.Range("E45") "9.17%"
Set-TextCell $ws.Range("G45") "13"

# Row 46
Set-TextCell $ws.Range("D46") "0.00006794"
Set-TextCell $ws.Range("E46") "2.70%"
Set-TextCell $ws.Range("G46") "13"

# Row 47
Set-TextCell $ws.Range("D47") "0.00000000750"
Set-TextCell $ws.Range("G47") "13"

# Row 48
Set-TextCell $ws.Range("D48") "0.003027"
Set-TextCell $ws.Range("E48") "5.33%"
Set-TextCell $ws.Range("G48") "13"

# Row 49
Set-TextCell $ws.Range("D49") "0.002280"
Set-TextCell $ws.Range("E49") "-5.07%"
Set-TextCell $ws.Range("G49") "13"

# Row 50
Set-TextCell $ws.Range("D50") "0.00002099"
Set-TextCell $ws.Range("G50") "13"

# Row 51
Set-TextCell $ws.Range("D51") "0.0001999"
Set-TextCell $ws.Range("G51") "13"
